$wb = $excel.ActiveWorkbook

# Update "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1674
$ws1.Range("F9").Value = 583

# Update "全部类型" sheet
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F3").Value = 1674
$ws2.Range("F9").Value = 583
